# Support multi-target upgrades with percent attack bonuses
# Add new unit rows to the Orc, NightElf, and Undead sheets.

$wb = $excel.ActiveWorkbook

# --- Orc sheet: add "Demolisher" as row 14 ---
$orc = $wb.Worksheets.Item("Orc")
$orcRow = 14
$orc.Cells.Item($orcRow, 1).Value = "Demolisher"
$orc.Cells.Item($orcRow, 2).Value = 220
$orc.Cells.Item($orcRow, 3).Value = 50
$orc.Cells.Item($orcRow, 4).Value = 4
$orc.Cells.Item($orcRow, 5).Value = 45
$orc.Cells.Item($orcRow, 6).Value = 72
$orc.Cells.Item($orcRow, 7).Value = 89
$orc.Cells.Item($orcRow, 8).Value = 425
$orc.Cells.Item($orcRow, 9).Value = 2
$orc.Cells.Item($orcRow, 10).Value = 600
$orc.Cells.Item($orcRow, 11).Value = "Siege"
$orc.Cells.Item($orcRow, 12).Value = 4.5
$orc.Cells.Item($orcRow, 13).Value = "Siege"
$orc.Cells.Item($orcRow, 14).Value = "Heavy"

# --- NightElf sheet: add "Mountain Giant" as row 13 ---
$nightElf = $wb.Worksheets.Item("NightElf")
$neRow = 13
$nightElf.Cells.Item($neRow, 1).Value = "Mountain Giant"
$nightElf.Cells.Item($neRow, 2).Value = 350
$nightElf.Cells.Item($neRow, 3).Value = 100
$nightElf.Cells.Item($neRow, 4).Value = 7
$nightElf.Cells.Item($neRow, 5).Value = 60
$nightElf.Cells.Item($neRow, 6).Value = 28
$nightElf.Cells.Item($neRow, 7).Value = 40
$nightElf.Cells.Item($neRow, 8).Value = 1600
$nightElf.Cells.Item($neRow, 9).Value = 6
$nightElf.Cells.Item($neRow, 10).Value = 100
$nightElf.Cells.Item($neRow, 11).Value = "Melee"
$nightElf.Cells.Item($neRow, 12).Value = 2.5
$nightElf.Cells.Item($neRow, 13).Value = "Normal"
$nightElf.Cells.Item($neRow, 14).Value = "Medium"

# --- Undead sheet: add "Obsidian Statue" as row 12 ---
$undead = $wb.Worksheets.Item("Undead")
$unRow = 12
$undead.Cells.Item($unRow, 1).Value = "Obsidian Statue"
$undead.Cells.Item($unRow, 2).Value = 200
$undead.Cells.Item($unRow, 3).Value = 35
$undead.Cells.Item($unRow, 4).Value = 3
$undead.Cells.Item($unRow, 5).Value = 30
$undead.Cells.Item($unRow, 6).Value = 7
$undead.Cells.Item($unRow, 7).Value = 8
$undead.Cells.Item($unRow, 8).Value = 500
$undead.Cells.Item($unRow, 9).Value = 4
$undead.Cells.Item($unRow, 10).Value = 575
$undead.Cells.Item($unRow, 11).Value = "Ranged"
$undead.Cells.Item($unRow, 12).Value = 2.1
$undead.Cells.Item($unRow, 13).Value = "Magic"
$undead.Cells.Item($unRow, 14).Value = "Heavy"

Write-Output "Added Demolisher, Mountain Giant, Obsidian Statue rows"
